$wb = $excel.ActiveWorkbook

# Rename the worksheets
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

$ws1.Name = "test1.txt_1"
$ws2.Name = "test2.txt_2"
$ws3.Name = "other_file_name.txt_3"

# --- Sheet1 (test1.txt_1) ---
# Clear previous content in row 1 columns beyond E, and set new values
$ws1.Cells.Clear()
$ws1.Range("A1").Value = "this"
$ws1.Range("B1").Value = "is"
$ws1.Range("C1").Value = "the"
$ws1.Range("D1").Value = "first"
$ws1.Range("E1").Value = "file"
$ws1.Range("A2").Value = "test1.txt"

# --- Sheet2 (test2.txt_2) ---
$ws2.Cells.Clear()
$ws2.Range("A1").Value = "second"
$ws2.Range("A2").Value = "file"
$ws2.Range("A3").Value = "test2.txt"
$ws2.Range("A4").Value = "this"
$ws2.Range("B4").Value = "is"
$ws2.Range("C4").Value = "the"
$ws2.Range("D4").Value = "second"
$ws2.Range("E4").Value = "file"
$ws2.Range("A5").Value = "the"
$ws2.Range("B5").Value = "second"
$ws2.Range("C5").Value = "one"

# --- Sheet3 (other_file_name.txt_3) ---
$ws3.Cells.Clear()
$ws3.Range("A1").Value = "the"
$ws3.Range("B1").Value = "file"
$ws3.Range("C1").Value = "can"
$ws3.Range("D1").Value = "have"
$ws3.Range("E1").Value = "other"
$ws3.Range("F1").Value = "file"
$ws3.Range("G1").Value = "names"
$ws3.Range("A2").Value = "it"
$ws3.Range("B2").Value = "has"
$ws3.Range("C2").Value = "to"
$ws3.Range("D2").Value = "end"
$ws3.Range("E2").Value = "with"
$ws3.Range("F2").Value = ".txt"
$ws3.Range("A3").Value = "(that"
$ws3.Range("B3").Value = "can"
$ws3.Range("C3").Value = "also"
$ws3.Range("D3").Value = "change)"

$wb.Save()
